{"js": "// Small wording fixes to \"Additional file 1\" (IH 24NOV2022).\n// Each change is a targeted find-and-replace of text within a single run,\n// using Word's search so surrounding run formatting (rFonts/sz/rtl/etc.) is\n// left untouched.\n\nconst body = context.document.body;\n\nconst replacements = [\n  {\n    find: \"the variant may be not be called\",\n    replace: \"the variant may not be called\",\n  },\n  {\n    find: \"primer design or library preparation chemistry.\",\n    replace: \"primer design or library preparation method.\",\n  },\n  {\n    find: \"fifteen tiles (RACE-like chemistry)\",\n    replace: \"fifteen tiles (RACE-like method)\",\n  },\n  {\n    find: \"67 million read pairs (RACE-like chemistry,\",\n    replace: \"67 million read pairs (RACE-like method,\",\n  },\n  {\n    find: \"We opted to use umi_tools for UMI grouping and satmut_utils for consensus deduplication of grouped reads using our flexible Python implementation.\",\n    replace: \"We opted to use umi_tools for UMI grouping and satmut_utils as a flexible implementation for consensus deduplication of grouped reads.\",\n  },\n  {\n    find: \"the user must manually filter the encompassing variant call\",\n    replace: \"the user may manually filter the encompassing variant call\",\n  },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for: ${find}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Small wording fixes to \"Additional file 1\" (IH 24NOV2022).\n# Each change is a targeted Find & Replace over the document body, leaving\n# surrounding run formatting untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"the variant may be not be called\"; Replace = \"the variant may not be called\" },\n    @{ Find = \"primer design or library preparation chemistry.\"; Replace = \"primer design or library preparation method.\" },\n    @{ Find = \"fifteen tiles (RACE-like chemistry)\"; Replace = \"fifteen tiles (RACE-like method)\" },\n    @{ Find = \"67 million read pairs (RACE-like chemistry,\"; Replace = \"67 million read pairs (RACE-like method,\" },\n    @{ Find = \"We opted to use umi_tools for UMI grouping and satmut_utils for consensus deduplication of grouped reads using our flexible Python implementation.\"; Replace = \"We opted to use umi_tools for UMI grouping and satmut_utils as a flexible implementation for consensus deduplication of grouped reads.\" },\n    @{ Find = \"the user must manually filter the encompassing variant call\"; Replace = \"the user may manually filter the encompassing variant call\" }\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($r.Find, $true, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2)\n    if (-not $found) {\n        throw \"No match found for: $($r.Find)\"\n    }\n}\n"}
